$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.798.56"
$ws.Range("D3").Value = "3.313.04"
$ws.Range("E3").Value = "  +6.25%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'599.28"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "'142.79"
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.311.42"
$ws.Range("E8").Value = "  +6.47%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").Value = "'5.56"
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("D12").Value = "'0.474"
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'34.78"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "3.856.50"
$ws.Range("E15").Value = "  +6.23%  "
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "3.315.80"
$ws.Range("E17").Value = "  +6.41%  "
$ws.Range("D18").Value = "63.886.20"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "'480.89"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "'0.735"
$ws.Range("E22").Value = "  +5.79%  "
$ws.Range("D23").Value = "'8.00"
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("D24").Value = "'13.56"
$ws.Range("E24").Value = "  +5.50%  "
$ws.Range("D25").Value = "'84.74"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("E28").Value = "  +5.28%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'8.16"
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("D31").Value = "'2.17"
$ws.Range("E31").Value = "  +4.75%  "
$ws.Range("D32").Value = "'29.31"
$ws.Range("E32").Value = "  +9.86%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +3.35%  "
$ws.Range("D37").Value = "0.0₃0754"
$ws.Range("E37").Value = "  +8.07%  "
$ws.Range("D38").Value = "'52.83"
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("D39").Value = "'0.0403"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("D40").Value = "'430.02"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").Value = "3.051.83"
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("E42").Value = "  +2.99%  "
$ws.Range("D43").Value = "'2.75"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "'0.266"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'2.21"
$ws.Range("E46").Value = "  +4.62%  "
$ws.Range("D47").Value = "'26.47"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").Value = "'36.19"
$ws.Range("E48").Value = "  +15.40%  "
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "  +2.36%  "
